# Upload Leave Card 12/27/2023 4:01 PM
# Leave-card workbook update for ROMILLA, EDITHA:
#  - Fill in EARNED (1.25) for Jul-Oct 2023 rows (111-114)
#  - Insert two new leave entries (SP(1-0-0), SL(2-0-0), SP(1-0-0)) spanning
#    Nov 2023, which pushes two brand-new blank rows into the monthly table
#  - Table1 / dimension grow from 145 to 147 rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$gFormula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# ---------------------------------------------------------------------
# 1) Fill EARNED (column C) for the Jul, Aug, Sep, Oct 2023 rows.
#    The mirrored "EARNED " column (G) is a calculated table column and
#    will recompute on its own once C holds a value.
# ---------------------------------------------------------------------
$ws.Range("C111").Value = 1.25
$ws.Range("C112").Value = 1.25
$ws.Range("C113").Value = 1.25
$ws.Range("C114").Value = 1.25

# ---------------------------------------------------------------------
# 2) Insert two new rows right after row 115 (i.e. before the old row
#    116), shifting every later row down by two. This is where the two
#    extra SL/SP leave rows for Nov 2023 will live.
# ---------------------------------------------------------------------
$ws.Range("A116:A117").EntireRow.Insert()

# The freshly-inserted rows 116:117 pick up a blended style from the
# insert; restore the same plain styling used by the rest of the table
# by pasting formats from row 118 (a normal, unmodified monthly row).
$ws.Range("A118:K118").Copy()
$ws.Range("A116:K117").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the calculated-column formula for EARNED (col G) on the two
# new rows (PasteSpecial of formats only does not carry formulas).
$ws.Range("G116").Formula = $gFormula
$ws.Range("G117").Formula = $gFormula

# ---------------------------------------------------------------------
# 3) Grow Table1 to cover the two new rows (A8:K145 -> A8:K147).
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K147"))

# Resizing the table while rows 146:147 (formerly 144:145, the last two
# rows of the sheet) sit just past the old boundary leaves their
# EARNED-mirror formula referencing the table in a stale way, producing
# a #VALUE!/#ERROR! result. Re-enter the same formula so it binds to
# the table correctly again now that it is an interior row.
$ws.Range("G146").Formula = $gFormula
$ws.Range("G147").Formula = $gFormula

# ---------------------------------------------------------------------
# 4) Give the REMARKS cells (K115:K117) that will hold dates the same
#    date-number-format style already used elsewhere in this column
#    (e.g. K93), then fill in the three new leave-card entries.
# ---------------------------------------------------------------------
$ws.Range("K93").Copy()
$ws.Range("K115:K117").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 115 (Nov 1, 2023 period) - Special Privilege leave filed 11/6/2023
$ws.Range("B115").Value = "SP(1-0-0)"
$ws.Range("K115").Value = "11/6/2023"

# Row 116 (new) - Sick Leave, 2 days, filed 11/17 & 11/20/2023
$ws.Range("B116").Value = "SL(2-0-0)"
$ws.Range("H116").Value = 2
$ws.Range("K116").Value = "11/17,20/2023"

# Row 117 (new) - Special Privilege leave filed 11/25/2023
$ws.Range("B117").Value = "SP(1-0-0)"
$ws.Range("K117").Value = "11/25/2023"

# ---------------------------------------------------------------------
# 5) Leave the cursor where the author last left it.
# ---------------------------------------------------------------------
$ws.Range("K117").Select()
